$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 36 (row 45): DC vs RR
$ws.Range("E45").Value = 60
$ws.Range("H45").Value = 100
$ws.Range("K45").Value = 40
$ws.Range("N45").Value = 20
$ws.Range("Q45").Value = 0
$ws.Range("T45").Value = 80

# Contest 37 (row 46): SRH vs PBKS
$ws.Range("E46").Value = 100
$ws.Range("H46").Value = 40
$ws.Range("K46").Value = 60
$ws.Range("N46").Value = 20
$ws.Range("Q46").Value = 0
$ws.Range("T46").Value = 80
